$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the claim's policy number and incident date (F2, H2) first so the
# shared-string table appends these new values ahead of the Ambiente/URL
# strings below (matches the order Excel produced when re-writing sharedStrings.xml).
# Use a leading apostrophe via .Formula so the text-quote style (s="3"/s="1")
# on these cells is preserved instead of Excel re-evaluating them as numbers/dates.
$ws.Range("F2").Formula = "'04104016408"
$ws.Range("H2").Formula = "'30/04/2021"

# Drop the "i-" preprod prefix from the Ambiente/URL columns for row 2.
$ws.Range("B2").Formula = "preproducciongestion.segurossura.com.ar"
$ws.Range("C2").Formula = "https://preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do"

# Add a hyperlink on C2 pointing at the new URL, matching the other rows (C3:C9).
$c2Style = $ws.Range("C2").Style
[void]$ws.Hyperlinks.Add($ws.Range("C2"), "https://preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do")
$ws.Range("C2").Style = $c2Style

# Reflect the user's final selection on the sheet.
[void]$ws.Range("C3").Select()
